$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(5, 8).Value = 234.7
$ws.Cells.Item(5, 9).Value = 144.5
$ws.Cells.Item(5, 11).Value = 144.5
$ws.Cells.Item(5, 13).Value = -29.5
$ws.Cells.Item(34, 8).Value = 9355
$ws.Cells.Item(34, 9).Value = 9355
$ws.Cells.Item(34, 11).Value = 9355
$ws.Cells.Item(34, 13).Value = -9152
$ws.Cells.Item(36, 8).Value = 9355
$ws.Cells.Item(36, 9).Value = 9355
$ws.Cells.Item(36, 11).Value = 9355
$ws.Cells.Item(36, 13).Value = -8640
$ws.Cells.Item(76, 8).Value = 6256
$ws.Cells.Item(76, 9).Value = 3650
$ws.Cells.Item(76, 10).Value = 7993.3335
$ws.Cells.Item(76, 11).Value = 3650
$ws.Cells.Item(76, 12).Value = 7993.3335
$ws.Cells.Item(76, 13).Value = -3335
$ws.Cells.Item(76, 14).Value = -8623.333500000001
$ws.Cells.Item(79, 8).Value = 6256
$ws.Cells.Item(79, 9).Value = 3650
$ws.Cells.Item(79, 10).Value = 7993.3335
$ws.Cells.Item(79, 11).Value = 3650
$ws.Cells.Item(79, 12).Value = 7993.3335
$ws.Cells.Item(79, 13).Value = -2558
$ws.Cells.Item(79, 14).Value = -10177.3335
$ws.Cells.Item(100, 8).Value = 5232.048
$ws.Cells.Item(100, 9).Value = 1688.5
$ws.Cells.Item(100, 10).Value = 8453.454
$ws.Cells.Item(100, 11).Value = 1688.5
$ws.Cells.Item(100, 12).Value = 8453.454
$ws.Cells.Item(100, 13).Value = -1147.5
$ws.Cells.Item(100, 14).Value = -9535.454
$ws.Cells.Item(107, 8).Value = 1215.6875
$ws.Cells.Item(107, 9).Value = 819.0769
$ws.Cells.Item(107, 11).Value = 819.0769
$ws.Cells.Item(107, 13).Value = 1100.9231
$ws.Cells.Item(135, 8).Value = 2090.2246
$ws.Cells.Item(135, 9).Value = 1363.9474
$ws.Cells.Item(135, 11).Value = 12275.5266
$ws.Cells.Item(135, 13).Value = -9740.526600000001
$ws.Cells.Item(137, 8).Value = 1497.1052
$ws.Cells.Item(137, 9).Value = 1291.1765
$ws.Cells.Item(137, 11).Value = 3873.5295
$ws.Cells.Item(137, 13).Value = -1323.5295

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 4051.1162
$ws.Cells.Item(32, 9).Value = 4089.475
$ws.Cells.Item(32, 10).Value = 3539.6667
$ws.Cells.Item(32, 11).Value = 4089.475
$ws.Cells.Item(32, 12).Value = 3539.6667
$ws.Cells.Item(32, 13).Value = -3802.475
$ws.Cells.Item(32, 14).Value = -4113.6667
$ws.Cells.Item(132, 8).Value = 6860
$ws.Cells.Item(132, 9).Value = 6860
$ws.Cells.Item(132, 11).Value = 20580
$ws.Cells.Item(132, 13).Value = -18050

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 1322.2142
$ws.Cells.Item(20, 9).Value = 1321.9
$ws.Cells.Item(20, 11).Value = 1321.9
$ws.Cells.Item(20, 13).Value = -1074.9
$ws.Cells.Item(94, 8).Value = 3278.4167
$ws.Cells.Item(94, 9).Value = 3540.0908
$ws.Cells.Item(94, 11).Value = 3540.0908
$ws.Cells.Item(94, 13).Value = -3089.0908

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(18, 8).Value = 99999
$ws.Cells.Item(18, 10).Value = 99999
$ws.Cells.Item(18, 12).Value = 99999
$ws.Cells.Item(18, 14).Value = -100459
$ws.Cells.Item(58, 8).Value = 2134.6667
$ws.Cells.Item(58, 9).Value = 1884.1428
$ws.Cells.Item(58, 11).Value = 1884.1428
$ws.Cells.Item(58, 13).Value = -1681.1428
$ws.Cells.Item(70, 8).Value = 99999
$ws.Cells.Item(70, 10).Value = 99999
$ws.Cells.Item(70, 12).Value = 99999
$ws.Cells.Item(70, 14).Value = -100629
$ws.Cells.Item(73, 8).Value = 99999
$ws.Cells.Item(73, 10).Value = 99999
$ws.Cells.Item(73, 12).Value = 99999
$ws.Cells.Item(73, 14).Value = -102183
$ws.Cells.Item(82, 8).Value = 67896.42999999999
$ws.Cells.Item(82, 10).Value = 70045.836
$ws.Cells.Item(82, 12).Value = 70045.836
$ws.Cells.Item(82, 14).Value = -70767.836
$ws.Cells.Item(85, 8).Value = 67896.42999999999
$ws.Cells.Item(85, 10).Value = 70045.836
$ws.Cells.Item(85, 12).Value = 70045.836
$ws.Cells.Item(85, 14).Value = -72541.836
$ws.Cells.Item(103, 8).Value = 64599.25
$ws.Cells.Item(103, 9).Value = 44199.5
$ws.Cells.Item(103, 11).Value = 44199.5
$ws.Cells.Item(103, 13).Value = -43027.5
$ws.Cells.Item(119, 8).Value = 77929
$ws.Cells.Item(119, 10).Value = 77929
$ws.Cells.Item(119, 12).Value = 77929
$ws.Cells.Item(119, 14).Value = -87605
$ws.Cells.Item(122, 8).Value = 2455.8125
$ws.Cells.Item(122, 9).Value = 1953.579
$ws.Cells.Item(122, 11).Value = 5860.737
$ws.Cells.Item(122, 13).Value = -3410.737
$ws.Cells.Item(136, 8).Value = 2134.6667
$ws.Cells.Item(136, 9).Value = 1884.1428
$ws.Cells.Item(136, 11).Value = 5652.428400000001
$ws.Cells.Item(136, 13).Value = -3102.428400000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 1458.1428
$ws.Cells.Item(5, 10).Value = 1389.25
$ws.Cells.Item(5, 12).Value = 4167.75
$ws.Cells.Item(5, 14).Value = -4391.75
$ws.Cells.Item(14, 8).Value = 10062.917
$ws.Cells.Item(14, 9).Value = 10062.917
$ws.Cells.Item(14, 11).Value = 30188.751
$ws.Cells.Item(14, 13).Value = -30015.751
$ws.Cells.Item(38, 8).Value = 686.6818
$ws.Cells.Item(38, 10).Value = 663.44446
$ws.Cells.Item(38, 12).Value = 1990.33338
$ws.Cells.Item(38, 14).Value = -2684.33338
$ws.Cells.Item(68, 8).Value = 1784.6
$ws.Cells.Item(68, 10).Value = 1906.5714
$ws.Cells.Item(68, 12).Value = 5719.7142
$ws.Cells.Item(68, 14).Value = -7341.7142
$ws.Cells.Item(71, 8).Value = 1784.6
$ws.Cells.Item(71, 10).Value = 1906.5714
$ws.Cells.Item(71, 12).Value = 17159.1426
$ws.Cells.Item(71, 14).Value = -25271.1426
$ws.Cells.Item(107, 8).Value = 2278649.5
$ws.Cells.Item(107, 10).Value = 3253993.5
$ws.Cells.Item(107, 12).Value = 9761980.5
$ws.Cells.Item(107, 14).Value = -9765820.5
$ws.Cells.Item(113, 8).Value = 1430.875
$ws.Cells.Item(113, 9).Value = 400
$ws.Cells.Item(113, 10).Value = 1578.1428
$ws.Cells.Item(113, 11).Value = 1200
$ws.Cells.Item(113, 12).Value = 4734.428400000001
$ws.Cells.Item(113, 13).Value = 970
$ws.Cells.Item(113, 14).Value = -9074.428400000001
$ws.Cells.Item(118, 8).Value = 6691.3
$ws.Cells.Item(118, 9).Value = 3731.111
$ws.Cells.Item(118, 11).Value = 11193.333
$ws.Cells.Item(118, 13).Value = -9950.332999999999
$ws.Cells.Item(122, 8).Value = 33332.1
$ws.Cells.Item(122, 9).Value = 66133
$ws.Cells.Item(122, 10).Value = 531.2
$ws.Cells.Item(122, 11).Value = 595197
$ws.Cells.Item(122, 12).Value = 4780.8
$ws.Cells.Item(122, 13).Value = -592747
$ws.Cells.Item(122, 14).Value = -9680.799999999999
$ws.Cells.Item(135, 8).Value = 1458.1428
$ws.Cells.Item(135, 10).Value = 1389.25
$ws.Cells.Item(135, 12).Value = 12503.25
$ws.Cells.Item(135, 14).Value = -17573.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(97, 8).Value = 1026.5385
$ws.Cells.Item(97, 9).Value = 619.3333
$ws.Cells.Item(97, 10).Value = 1942.75
$ws.Cells.Item(97, 11).Value = 619.3333
$ws.Cells.Item(97, 12).Value = 1942.75
$ws.Cells.Item(97, 13).Value = -123.3333
$ws.Cells.Item(97, 14).Value = -2934.75
$ws.Cells.Item(119, 8).Value = 100380
$ws.Cells.Item(119, 10).Value = 100380
$ws.Cells.Item(119, 12).Value = 100380
$ws.Cells.Item(119, 14).Value = -110056

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(46, 8).Value = 1397.1666
$ws.Cells.Item(46, 10).Value = 1398
$ws.Cells.Item(46, 12).Value = 1398
$ws.Cells.Item(46, 14).Value = -1774
$ws.Cells.Item(82, 8).Value = 7577.4546
$ws.Cells.Item(82, 9).Value = 4878.4287
$ws.Cells.Item(82, 11).Value = 4878.4287
$ws.Cells.Item(82, 13).Value = -4517.4287
$ws.Cells.Item(85, 8).Value = 7577.4546
$ws.Cells.Item(85, 9).Value = 4878.4287
$ws.Cells.Item(85, 11).Value = 4878.4287
$ws.Cells.Item(85, 13).Value = -3630.4287
$ws.Cells.Item(93, 8).Value = 2419387.5
$ws.Cells.Item(93, 9).Value = 2097.2222
$ws.Cells.Item(93, 10).Value = 3973359.8
$ws.Cells.Item(93, 11).Value = 2097.2222
$ws.Cells.Item(93, 12).Value = 3973359.8
$ws.Cells.Item(93, 13).Value = -849.2222000000002
$ws.Cells.Item(93, 14).Value = -3975855.8
$ws.Cells.Item(100, 8).Value = 13909454
$ws.Cells.Item(100, 9).Value = 7268.0835
$ws.Cells.Item(100, 11).Value = 7268.0835
$ws.Cells.Item(100, 13).Value = -6727.0835
$ws.Cells.Item(132, 8).Value = 2810.8572
$ws.Cells.Item(132, 9).Value = 2266.2
$ws.Cells.Item(132, 10).Value = 4172.5
$ws.Cells.Item(132, 11).Value = 6798.599999999999
$ws.Cells.Item(132, 12).Value = 12517.5
$ws.Cells.Item(132, 13).Value = -4268.599999999999
$ws.Cells.Item(132, 14).Value = -17577.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(81, 8).Value = 1281
$ws.Cells.Item(81, 9).Value = 1421.5
$ws.Cells.Item(81, 11).Value = 2843
$ws.Cells.Item(81, 13).Value = -1782
$ws.Cells.Item(84, 8).Value = 1281
$ws.Cells.Item(84, 9).Value = 1421.5
$ws.Cells.Item(84, 11).Value = 14215
$ws.Cells.Item(84, 13).Value = -8911
$ws.Cells.Item(100, 8).Value = 704.7778
$ws.Cells.Item(100, 9).Value = 732.3333
$ws.Cells.Item(100, 10).Value = 649.6667
$ws.Cells.Item(100, 11).Value = 1464.6666
$ws.Cells.Item(100, 12).Value = 1299.3334
$ws.Cells.Item(100, 13).Value = -923.6666
$ws.Cells.Item(100, 14).Value = -2381.3334
$ws.Cells.Item(113, 8).Value = 1409.8235
$ws.Cells.Item(113, 9).Value = 1238.5834
$ws.Cells.Item(113, 11).Value = 3715.7502
$ws.Cells.Item(113, 13).Value = -1545.7502
$ws.Cells.Item(122, 8).Value = 2042.5834
$ws.Cells.Item(122, 9).Value = 1815.619
$ws.Cells.Item(122, 11).Value = 5446.857
$ws.Cells.Item(122, 13).Value = -2996.857
$ws.Cells.Item(126, 8).Value = 4047.238
$ws.Cells.Item(126, 9).Value = 3429.6428
$ws.Cells.Item(126, 10).Value = 5282.4287
$ws.Cells.Item(126, 11).Value = 10288.9284
$ws.Cells.Item(126, 12).Value = 15847.2861
$ws.Cells.Item(126, 13).Value = -7818.928400000001
$ws.Cells.Item(126, 14).Value = -20787.2861
